$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "UAT12"
$ws.Range("C14").Value = "RRD"
$ws.Range("D14").Value = "PROD18"
$ws.Range("E14").Value = "JRD"
$ws.Range("F14").Value = "Trades"
$ws.Range("G14").Value = "Partition Copy"
$ws.Range("H14").Value = "10-05-2020 22:29:18"
$ws.Range("I14").Value = "10-05-2020 22:29:18"
$ws.Range("J14").Value = "In Progress"
$ws.Range("K14").Value = "Export in Progress"
$ws.Range("L14").Value = "Import in Progress"
